$wb = $excel.ActiveWorkbook

# Fix typo "uniquie" -> "unique" in the taxon_id definition on the glossary sheet
$glossary = $wb.Worksheets.Item("glossary")
$glossary.Range("B49").Value = "The unique identifier for a single taxon. This column should include ALL unique taxon_id entries from the fish survey data sheet"

# Remove the duplicate/incorrect taxon_id row (row 50) which held an
# out-of-date definition ("A unique string (typically 3 digits)...").
# Deleting the row shifts the transect row (51) up to 50 and the blank
# trailing row (52) up to 51, matching the target layout.
$glossary.Rows(50).Delete()
